$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 970
$ws.Range("I88").Value = 862.5
$ws.Range("J88").Value = 996.875
$ws.Range("K88").Value = 862.5
$ws.Range("L88").Value = 996.875
$ws.Range("M88").Value = -456.5
$ws.Range("N88").Value = -1808.875
$ws.Range("H91").Value = 970
$ws.Range("I91").Value = 862.5
$ws.Range("J91").Value = 996.875
$ws.Range("K91").Value = 862.5
$ws.Range("L91").Value = 996.875
$ws.Range("M91").Value = 541.5
$ws.Range("N91").Value = -3804.875
$ws.Range("H128").Value = 15797.5
$ws.Range("J128").Value = 15797.5
$ws.Range("L128").Value = 15797.5
$ws.Range("N128").Value = -25757.5
$ws.Range("H129").Value = 937.68604
$ws.Range("J129").Value = 1011.8947
$ws.Range("L129").Value = 3035.6841
$ws.Range("N129").Value = -13035.6841
$ws.Range("H137").Value = 1278.25
$ws.Range("I137").Value = 995.561
$ws.Range("J137").Value = 2934
$ws.Range("K137").Value = 2986.683
$ws.Range("L137").Value = 8802
$ws.Range("M137").Value = -436.683
$ws.Range("N137").Value = -13902

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3251.1667
$ws.Range("J88").Value = 3251.1667
$ws.Range("L88").Value = 3251.1667
$ws.Range("N88").Value = -4063.1667
$ws.Range("H91").Value = 3251.1667
$ws.Range("J91").Value = 3251.1667
$ws.Range("L91").Value = 3251.1667
$ws.Range("N91").Value = -6059.1667
$ws.Range("H102").Value = 2569.111
$ws.Range("I102").Value = 2569.111
$ws.Range("K102").Value = 2569.111
$ws.Range("M102").Value = -947.1109999999999
$ws.Range("H132").Value = 1966.9844
$ws.Range("I132").Value = 917.7857
$ws.Range("J132").Value = 3970
$ws.Range("K132").Value = 2753.3571
$ws.Range("L132").Value = 11910
$ws.Range("M132").Value = -223.3571000000002
$ws.Range("N132").Value = -16970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2150.5625
$ws.Range("I86").Value = 2283.5
$ws.Range("J86").Value = 1751.75
$ws.Range("K86").Value = 2283.5
$ws.Range("L86").Value = 1751.75
$ws.Range("M86").Value = -1160.5
$ws.Range("N86").Value = -3997.75
$ws.Range("H89").Value = 2150.5625
$ws.Range("I89").Value = 2283.5
$ws.Range("J89").Value = 1751.75
$ws.Range("K89").Value = 11417.5
$ws.Range("L89").Value = 8758.75
$ws.Range("M89").Value = -5801.5
$ws.Range("N89").Value = -19990.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2989.3235
$ws.Range("I31").Value = 2714.1052
$ws.Range("K31").Value = 2714.1052
$ws.Range("M31").Value = -2419.1052
$ws.Range("H34").Value = 2989.3235
$ws.Range("I34").Value = 2714.1052
$ws.Range("K34").Value = 2714.1052
$ws.Range("M34").Value = -2512.1052
$ws.Range("H58").Value = 1742.6459
$ws.Range("I58").Value = 1017.8
$ws.Range("J58").Value = 2530.5217
$ws.Range("K58").Value = 1017.8
$ws.Range("L58").Value = 2530.5217
$ws.Range("M58").Value = -814.8
$ws.Range("N58").Value = -2936.5217
$ws.Range("H62").Value = 5542.857
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 6760
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 6760
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -8008
$ws.Range("H65").Value = 5542.857
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 6760
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 33800
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -40040
$ws.Range("H134").Value = 2933.45
$ws.Range("I134").Value = 2772.1875
$ws.Range("K134").Value = 8316.5625
$ws.Range("M134").Value = -5781.5625
$ws.Range("H136").Value = 1742.6459
$ws.Range("I136").Value = 1017.8
$ws.Range("J136").Value = 2530.5217
$ws.Range("K136").Value = 3053.4
$ws.Range("L136").Value = 7591.5651
$ws.Range("M136").Value = -503.3999999999996
$ws.Range("N136").Value = -12691.5651

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 603.8182
$ws.Range("I122").Value = 330.375
$ws.Range("K122").Value = 2973.375
$ws.Range("M122").Value = -523.375
$ws.Range("H131").Value = 2625.2388
$ws.Range("I131").Value = 285
$ws.Range("J131").Value = 2697.246
$ws.Range("K131").Value = 855
$ws.Range("L131").Value = 8091.738
$ws.Range("M131").Value = 4185
$ws.Range("N131").Value = -18171.738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24337.291
$ws.Range("I70").Value = 32286.467
$ws.Range("J70").Value = 11088.667
$ws.Range("K70").Value = 32286.467
$ws.Range("L70").Value = 11088.667
$ws.Range("M70").Value = -32016.467
$ws.Range("N70").Value = -11628.667
$ws.Range("H73").Value = 24337.291
$ws.Range("I73").Value = 32286.467
$ws.Range("J73").Value = 11088.667
$ws.Range("K73").Value = 32286.467
$ws.Range("L73").Value = 11088.667
$ws.Range("M73").Value = -31350.467
$ws.Range("N73").Value = -12960.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 69153.47
$ws.Range("I7").Value = 112922
$ws.Range("J7").Value = 3500.6667
$ws.Range("K7").Value = 112922
$ws.Range("L7").Value = 3500.6667
$ws.Range("M7").Value = -112810
$ws.Range("N7").Value = -3724.6667
$ws.Range("H126").Value = 69153.47
$ws.Range("I126").Value = 112922
$ws.Range("J126").Value = 3500.6667
$ws.Range("K126").Value = 338766
$ws.Range("L126").Value = 10502.0001
$ws.Range("M126").Value = -336296
$ws.Range("N126").Value = -15442.0001
$ws.Range("H136").Value = 2625.7715
$ws.Range("I136").Value = 3438.4119
$ws.Range("J136").Value = 1858.2778
$ws.Range("K136").Value = 10315.2357
$ws.Range("L136").Value = 5574.8334
$ws.Range("M136").Value = -7765.235700000001
$ws.Range("N136").Value = -10674.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1467.5
$ws.Range("I96").Value = 1435
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1435
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -62
$ws.Range("N96").Value = -4246
$ws.Range("H136").Value = 9262249
$ws.Range("I136").Value = 15873847
$ws.Range("K136").Value = 47621541
$ws.Range("M136").Value = -47618991
